# Reconciliation exceptions sheet update:
#   - new transaction rows appended from fuzzy-matched reconciliation run
#     (9/12, 12/12 x2, 14/12 extra rows) and "Monto" values refreshed for
#     several already-present rows (see commit message).
#   - Rewritten row-by-row against the target snapshot so shared-string /
#     row ordering matches exactly what Excel produced when the rows were
#     appended and re-sorted by date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writing a literal-looking-like-a-date or numeric-looking string via
# Range.Value triggers Excel's "smart" type coercion (dates become serials,
# a leading '+'/'-' etc. can become numbers) and also stamps a NumberFormat-
# derived style on the cell. Route text through a throw-away formula +
# paste-special-values instead so the stored cell is a plain shared string
# with no style override, exactly like the other text cells already in
# this sheet.
function Set-TextCell($cell, [string]$val) {
    $escaped = $val -replace '"', '""'
    $cell.Formula = '=""&"' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$rows = @(
    @(2, '02/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'LIQUIDACION MP 820442', '2318096.87', 'REF597889', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(3, '02/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'MERPAG*WO', '486161.57', 'REF326387', 'Sin coincidencia suficiente', 'Agregar alias a tabla parametrica'),
    @(4, '03/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'LIQUIDACION MP 757694', '1759015.55', 'REF255061', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(5, '04/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'LIQUIDACION MP 625600', '189293.73', 'REF389279', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(6, '04/12/2025', 'Banco Galicia', 'CREDITO', 'cobranza', 'TRANSF TERCEROS CBU 9792513271', '229655.18', 'REF685472', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(7, '05/12/2025', 'Banco Galicia', 'CREDITO', 'cobranza', 'TRANSF TERCEROS CBU 3331279757', '107405.99', 'REF650744', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(8, '06/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'LIQUIDACION MP 398433', '344093.64', 'REF472358', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(9, '07/12/2025', 'Banco Galicia', 'CREDITO', 'cobranza', 'TRANSF TERCEROS CBU 0868034764', '103746.76', 'REF857750', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(10, '07/12/2025', 'Banco Galicia', 'CREDITO', 'cobranza', 'TRANSF TERCEROS CBU 0678831712', '2039783.71', 'REF952227', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(11, '07/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'LIQUIDACION MP 889410', '330194.04', 'REF724041', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(12, '08/12/2025', 'Banco Santander', 'CREDITO', 'cobranza', 'TRANSF TERCEROS CBU 8463920867', '551573.16', 'REF112394', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(13, '08/12/2025', 'Banco Galicia', 'CREDITO', 'cobranza', 'TRANSF TERCEROS CBU 0049082404', '403899.15', 'REF139569', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(14, '09/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'MERPAG*DADDYBEBIDASCORDILLE', '340845.75', 'REF952117', 'Sin coincidencia suficiente', 'Agregar alias a tabla parametrica'),
    @(15, '10/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'LIQUIDACION MP 748314', '379163.43', 'REF656268', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(16, '12/12/2025', 'Banco Galicia', 'CREDITO', 'cobranza', 'TRANSF TERCEROS CBU 1465669380', '411936.75', 'REF576849', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(17, '12/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'MERPAG*EDUARDOBRUNOVMARIA', '704337.95', 'REF823255', 'Sin coincidencia suficiente', 'Agregar alias a tabla parametrica'),
    @(18, '12/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'LIQUIDACION MP 427841', '284955.85', 'REF678944', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(19, '10/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'LIQUIDACION MP 385037', '414032.01', 'REF976303', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(20, '13/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'MERPAG*PRI', '29038053.54', 'REF866437', 'Sin coincidencia suficiente', 'Agregar alias a tabla parametrica'),
    @(21, '14/12/2025', 'Banco Galicia', 'CREDITO', 'cobranza', 'TRANSF TERCEROS CBU 6008350583', '187821.77', 'REF275133', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(22, '14/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'MERPAG*LAREINAEMPANADAS', '381603.21', 'REF517754', 'Sin coincidencia suficiente', 'Agregar alias a tabla parametrica'),
    @(23, '16/12/2025', 'Banco Santander', 'CREDITO', 'cobranza', 'TRANSF TERCEROS CBU 2410386571', '683935.83', 'REF405844', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(24, '16/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'LIQUIDACION MP 819335', '368338.83', 'REF339669', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(25, '18/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'LIQUIDACION MP 765678', '13396930.09', 'REF796405', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(26, '18/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'LIQUIDACION MP 178086', '400349.49', 'REF296059', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(27, '25/12/2025', 'Banco Galicia', 'CREDITO', 'cobranza', 'TRANSF TERCEROS CBU 6765861793', '395757.53', 'REF493696', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(28, '27/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'LIQUIDACION MP 783112', '237743.08', 'REF413305', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(29, '30/12/2025', 'Mercado Pago', 'CREDITO', 'cobranza', 'LIQUIDACION MP 618927', '857955.72', 'REF504630', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(30, '30/12/2025', 'Banco Galicia', 'CREDITO', 'cobranza', 'TRANSF TERCEROS CBU 8065787598', '362381.89', 'REF655149', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(31, '30/12/2025', 'Banco Galicia', 'CREDITO', 'cobranza', 'TRANSF TERCEROS CBU 1744518705', '1567777.27', 'REF897423', 'Sin match de identidad', 'Agregar alias a tabla parametrica'),
    @(32, '30/12/2025', 'Banco Galicia', 'CREDITO', 'cobranza', 'TRANSF MO', '202141.43', 'REF398916', 'Sin coincidencia suficiente', 'Agregar alias a tabla parametrica')
)

foreach ($row in $rows) {
    $r = $row[0]
    Set-TextCell $ws.Cells.Item($r, 1) $row[1]   # Fecha
    Set-TextCell $ws.Cells.Item($r, 2) $row[2]   # Banco
    Set-TextCell $ws.Cells.Item($r, 3) $row[3]   # Tipo
    Set-TextCell $ws.Cells.Item($r, 4) $row[4]   # Clasificacion
    Set-TextCell $ws.Cells.Item($r, 5) $row[5]   # Descripcion Original
    $ws.Cells.Item($r, 6).Value = [double]$row[6] # Monto (numeric)
    Set-TextCell $ws.Cells.Item($r, 7) $row[7]   # Referencia
    Set-TextCell $ws.Cells.Item($r, 8) $row[8]   # Detalle
    Set-TextCell $ws.Cells.Item($r, 9) $row[9]   # Accion Sugerida
}
